# The sheet's first row changes from raw numeric/placeholder data into a
# text header row, while rows 2-4 keep their existing values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Cilveki"
$ws.Range("B1").Value = "Nauda"
$ws.Range("C1").Value = "Durvis"
$ws.Range("D1").Value = "Piena litri"

# Move the active selection to match the saved workbook state (F3).
$ws.Range("F3").Select() | Out-Null
